$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Add work hours & descriptions for the two previously-empty rows (25/7/2017 and 26/7/2017)
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = "Fixed some bugs with editing disabling and tab selection"

$ws.Range("B16").Value = 8
$ws.Range("C16").Value = "Added scrollable tree view & branch item prefab"

# Update the selected cell in the sheet view
$ws.Range("B17").Select()

# Recalculate so the Total: formula in B38 reflects the new values
$excel.Calculate()
